$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '35.144.23'
$ws.Range('E2').Value = '  +0.27%  '

Set-TextValue $ws.Range('D3') '1.852.61'
$ws.Range('E3').Value = '  +1.39%  '

$ws.Range('E4').Value = '  +0.45%  '

Set-TextValue $ws.Range('D5') '237.64'
$ws.Range('E5').Value = '  +2.68%  '

$ws.Range('E6').Value = '  +0.36%  '

$ws.Range('E7').Value = '  +0.45%  '

Set-TextValue $ws.Range('D8') '41.99'
$ws.Range('E8').Value = '  +4.22%  '

Set-TextValue $ws.Range('D9') '0.326'
$ws.Range('E9').Value = '  +0.87%  '

$ws.Range('E10').Value = '  +1.09%  '

Set-TextValue $ws.Range('D11') '0.0987'
$ws.Range('E11').Value = '  -0.56%  '

Set-TextValue $ws.Range('D12') '2.121.45'
$ws.Range('E12').Value = '  +1.50%  '

$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D13') '11.38'
$ws.Range('E13').Value = '  +0.16%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D14') '1.847.16'
$ws.Range('E14').Value = '  +1.09%  '

Set-TextValue $ws.Range('D15') '0.675'
$ws.Range('E15').Value = '  +0.72%  '

Set-TextValue $ws.Range('D16') '4.71'
$ws.Range('E16').Value = '  +0.80%  '

Set-TextValue $ws.Range('D17') '35.082.76'
$ws.Range('E17').Value = '  -0.05%  '

Set-TextValue $ws.Range('D18') '69.88'
$ws.Range('E18').Value = '  +0.13%  '

$ws.Range('E19').Value = '  +0.36%  '

Set-TextValue $ws.Range('D20') '240.66'
$ws.Range('E20').Value = '  -0.16%  '

Set-TextValue $ws.Range('D21') '12.22'
$ws.Range('E21').Value = '  +0.77%  '

Set-TextValue $ws.Range('D22') '4.72'
$ws.Range('E22').Value = '  -1.15%  '

$ws.Range('E23').Value = '  +0.43%  '

$ws.Range('E24').Value = '  +0.32%  '

Set-TextValue $ws.Range('D25') '168.32'
$ws.Range('E25').Value = '  -3.47%  '

$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D26') '1.84'
$ws.Range('E26').Value = '  +20.44%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D27') '7.97'
$ws.Range('E27').Value = '  +1.73%  '

Set-TextValue $ws.Range('D28') '17.60'
$ws.Range('E28').Value = '  +1.04%  '

$ws.Range('E29').Value = '  -1.22%  '

$ws.Range('E30').Value = '  +0.48%  '

Set-TextValue $ws.Range('D31') '0.0554'
$ws.Range('E31').Value = '  +0.20%  '

Set-TextValue $ws.Range('D32') '3.97'
$ws.Range('E32').Value = '  -0.86%  '

$ws.Range('E33').Value = '  +0.44%  '

Set-TextValue $ws.Range('D34') '1.78'
$ws.Range('E34').Value = '  +28.55%  '

$ws.Range('E35').Value = '  +16.82%  '

$ws.Range('E36').Value = '  +8.99%  '

Set-TextValue $ws.Range('D37') '1.30'
$ws.Range('E37').Value = '  +5.20%  '

$ws.Range('E38').Value = '  +7.20%  '

$ws.Range('B39').Value = 'Aave'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D39') '90.06'
$ws.Range('E39').Value = '  -3.81%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D40') '0.0200'
$ws.Range('E40').Value = '  +2.73%  '

Set-TextValue $ws.Range('D41') '1.341.60'
$ws.Range('E41').Value = '  -0.06%  '

Set-TextValue $ws.Range('D42') '14.81'
$ws.Range('E42').Value = '  -0.28%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D43') '2.30'
$ws.Range('E43').Value = '  +0.64%  '

$ws.Range('B44').Value = 'Gas'
$ws.Range('C44').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
Set-TextValue $ws.Range('D44') '12.76'
$ws.Range('E44').Value = '  +49.33%  '

$ws.Range('E45').Value = '  +0.19%  '

$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D46') '0.0555'
$ws.Range('E46').Value = '  +6.33%  '

$ws.Range('B47').Value = 'MXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D47') '2.74'
$ws.Range('E47').Value = '  -1.00%  '

Set-TextValue $ws.Range('D48') '6.46'
$ws.Range('E48').Value = '  +2.49%  '

Set-TextValue $ws.Range('D49') '2.034.50'
$ws.Range('E49').Value = '  +1.42%  '

$ws.Range('E50').Value = '  +1.05%  '

$ws.Range('E51').Value = '  +0.45%  '
